$d = $word.ActiveDocument

# The paragraph "▲表 8-2-8 廠商" has the "▲" in its own run, immediately
# followed by a separate run containing "表 8-2-8 廠商". We need to remove
# just that leading "▲" run/character, leaving "表 8-2-8 廠商" (and its
# original run/formatting) completely untouched. "▲表 8-2-8" is unique in
# the document, so locate it with Find, then shrink the found range down
# to only the leading "▲" character before deleting it (a plain
# Find/Replace would instead merge the two runs' text/formatting).
$r = $d.Content
$r.Find.Execute("▲表 8-2-8", $true, $false, $false, $false, $false,
                 $true, 1, $false, "", 0)
# $r now spans the found text "▲表 8-2-8" (8 characters); shrink it down
# to just the leading "▲" character and delete that, leaving the
# "表 8-2-8" (and its original run/formatting) untouched.
[void]$r.MoveEnd(1, -7)
$r.Delete()
